$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Per-cell updates derived from the commit diff: a routine refresh of the
# "Price" (D) / "Volume(1h)" (E) columns for most coins, plus one swap of
# row 33/34 identities (EthereumClassic <-> ImmutableX, incl. their D/E).
$changes = @{
    2 = @{ D = "68.885.86"; E = "  -0.43%  " }
    3 = @{ D = "3.860.11"; E = "  +2.54%  " }
    4 = @{ D = "1.00"; E = "  +0.00%  " }
    5 = @{ D = "601.09"; E = "  -0.29%  " }
    6 = @{ D = "162.10"; E = "  -3.32%  " }
    7 = @{ D = "3.864.55"; E = "  +2.71%  " }
    8 = @{ E = "  +0.01%  " }
    10 = @{ E = "  -1.23%  " }
    11 = @{ E = "  -1.58%  " }
    12 = @{ E = "  -0.12%  " }
    13 = @{ D = "36.94"; E = "  -3.05%  " }
    14 = @{ E = "  -1.93%  " }
    15 = @{ D = "4.507.33"; E = "  +2.67%  " }
    16 = @{ D = "3.850.87"; E = "  +2.49%  " }
    17 = @{ D = "69.055.27"; E = "  -0.25%  " }
    19 = @{ D = "11.50"; E = "  +1.87%  " }
    20 = @{ E = "  -0.26%  " }
    21 = @{ D = "17.11"; E = "  -1.22%  " }
    22 = @{ D = "484.64"; E = "  -2.01%  " }
    23 = @{ D = "0.720"; E = "  -1.37%  " }
    24 = @{ D = "0.0000160"; E = "  +4.61%  " }
    25 = @{ D = "83.97" }
    26 = @{ D = "2.25"; E = "  -2.76%  " }
    27 = @{ D = "12.10" }
    28 = @{ D = "1.00"; E = "  -0.01%  " }
    29 = @{ E = "  -1.34%  " }
    30 = @{ D = "2.96" }
    31 = @{ D = "7.94"; E = "  -2.38%  " }
    32 = @{ D = "4.006.54"; E = "  +2.57%  " }
    33 = @{ B = "EthereumClassic"; C = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"; D = "32.23"; E = "  +1.79%  " }
    34 = @{ B = "ImmutableX"; C = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"; D = "2.37"; E = "  -4.35%  " }
    35 = @{ D = "3.807.53"; E = "  +2.91%  " }
    36 = @{ E = "  -1.31%  " }
    37 = @{ E = "  +0.80%  " }
    38 = @{ E = "  +1.86%  " }
    39 = @{ E = "  -2.03%  " }
    40 = @{ E = "  -0.11%  " }
    41 = @{ D = "0.320"; E = "  -1.96%  " }
    42 = @{ D = "437.17"; E = "  +1.60%  " }
    43 = @{ E = "  -1.66%  " }
    44 = @{ D = "48.49"; E = "  -0.39%  " }
    45 = @{ E = "  -0.51%  " }
    47 = @{ D = "8.40"; E = "  -1.02%  " }
    48 = @{ D = "143.33"; E = "  +1.75%  " }
    49 = @{ D = "2.838.03"; E = "  +1.51%  " }
    50 = @{ E = "  +1.71%  " }
    51 = @{ D = "25.91"; E = "  +12.02%  " }
}

# Columns whose new text would otherwise be auto-coerced into a Number by
# Excel's type inference (e.g. "1.00" -> 1, "0.720" -> 0.72), which would
# silently drop the significant trailing/leading zeros the source data needs.
# Pre-formatting just those cells as Text keeps the literal string intact.
$textForceCells = @(
    "D4", "D5", "D6", "D13", "D19", "D21", "D22", "D23", "D24", "D25", "D26", "D27", "D28", "D30", "D31", "D33", "D34", "D41", "D42", "D44", "D47", "D48", "D51"
)
foreach ($cellRef in $textForceCells) {
    $ws.Range($cellRef).NumberFormat = "@"
}

foreach ($row in $changes.Keys) {
    $rowData = $changes[$row]
    foreach ($col in $rowData.Keys) {
        $ws.Range("$col$row").Value = $rowData[$col]
    }
}
